$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to remain text (matching the original inline-string cells)
# while we overwrite their values, then strip the temporary formatting back
# off so the cells keep their original (unstyled) look.
$ws.Range("A2:D4").NumberFormat = "@"

$ws.Range("B2").Value = "-2.0"
$ws.Range("C2").Value = "-8.0"
$ws.Range("D2").Value = "2.0003606680642"

$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "0.000360668064204"
$ws.Range("C3").Value = "4.6916225675279003e-11"
$ws.Range("D3").Value = "0.0003576882720103"

$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "2.97979219374231e-06"
$ws.Range("C4").Value = "2.64580561779787e-17"
$ws.Range("D4").Value = "nan"

$ws.Range("A2:D4").Style = "Normal"

# Remove the now-obsolete iteration rows 5-8 (table shrinks to A1:D4)
$ws.Rows("5:8").Delete()
